$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.993.55"
$ws.Range("E2").Value = "  +3.26%  "

$ws.Range("D3").Value = "2.341.65"
$ws.Range("E3").Value = "  +2.59%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").Value = "'313.11"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").Value = "'108.69"
$ws.Range("E6").Value = "  +3.35%  "

$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "'0.619"

$ws.Range("D10").Value = "'41.38"
$ws.Range("E10").Value = "  +4.52%  "

$ws.Range("D11").Value = "'0.0919"
$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").Value = "'8.60"
$ws.Range("E12").Value = "  +2.53%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  -1.09%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'1.01"
$ws.Range("E14").Value = "  +1.58%  "

$ws.Range("D15").Value = "'15.50"
$ws.Range("E15").Value = "  +2.09%  "

$ws.Range("D16").Value = "2.696.99"
$ws.Range("E16").Value = "  +2.56%  "

$ws.Range("D17").Value = "2.338.69"
$ws.Range("E17").Value = "  +3.01%  "

$ws.Range("D18").Value = "43.904.93"
$ws.Range("E18").Value = "  +3.17%  "

$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  +2.97%  "

$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").Value = "'12.97"
$ws.Range("E21").Value = "  -4.62%  "

$ws.Range("D22").Value = "'74.29"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  -1.13%  "

$ws.Range("D24").Value = "'268.81"
$ws.Range("E24").Value = "  +1.75%  "

$ws.Range("E25").Value = "  +3.78%  "

$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("D27").Value = "'7.55"
$ws.Range("E27").Value = "  +6.69%  "

$ws.Range("D28").Value = "'11.16"
$ws.Range("E28").Value = "  +3.40%  "

$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("D30").Value = "'38.87"
$ws.Range("E30").Value = "  +4.99%  "

$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").Value = "'168.30"
$ws.Range("E32").Value = "  +1.25%  "

$ws.Range("D33").Value = "'0.0887"
$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").Value = "'2.79"
$ws.Range("E34").Value = "  +8.01%  "

$ws.Range("E35").Value = "  +1.32%  "

$ws.Range("D36").Value = "'4.76"
$ws.Range("E36").Value = "  +5.21%  "

$ws.Range("E37").Value = "  -0.26%  "

$ws.Range("D38").Value = "'0.0365"
$ws.Range("E38").Value = "  +4.37%  "

$ws.Range("D39").Value = "'2.88"
$ws.Range("E39").Value = "  +8.68%  "

$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("E41").Value = "  +8.76%  "

$ws.Range("D42").Value = "'105.08"
$ws.Range("E42").Value = "  +10.67%  "

$ws.Range("D43").Value = "'0.238"
$ws.Range("E43").Value = "  +3.06%  "

$ws.Range("D44").Value = "'71.71"
$ws.Range("E44").Value = "  +1.71%  "

$ws.Range("D45").Value = "'13.29"
$ws.Range("E45").Value = "  +9.14%  "

$ws.Range("E46").Value = "  +0.54%  "

$ws.Range("D47").Value = "'113.94"
$ws.Range("E47").Value = "  +0.40%  "

$ws.Range("D48").Value = "1.661.96"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("D49").Value = "'8.99"
$ws.Range("E49").Value = "  +2.96%  "

$ws.Range("D50").Value = "'76.50"
$ws.Range("E50").Value = "  -3.36%  "

$ws.Range("D51").Value = "'0.216"
$ws.Range("E51").Value = "  +14.06%  "
